$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows after row 300, pushing rows 301-350 down to 303-352
$ws.Rows.Item(301).Resize(2).Insert()

# Row 301 gets a copy of the data currently in row 299 (the original, not-yet-edited values)
# Row 302 gets a copy of the data currently in row 300
for ($col = 1; $col -le 18; $col++) {
    $ws.Cells.Item(301, $col).Value2 = $ws.Cells.Item(299, $col).Value2
    $ws.Cells.Item(302, $col).Value2 = $ws.Cells.Item(300, $col).Value2
}

# Now update row 299 with its new values
$ws.Cells.Item(299, 4).Value2 = 44476
$ws.Cells.Item(299, 10).Value2 = 600
$ws.Cells.Item(299, 11).Value2 = 16000
$ws.Cells.Item(299, 12).Value2 = 17000
$ws.Cells.Item(299, 13).Value2 = 16500
$ws.Cells.Item(299, 16).Value2 = 275

# Update row 300 with its new values
$ws.Cells.Item(300, 4).Value2 = 44476
$ws.Cells.Item(300, 10).Value2 = 360
$ws.Cells.Item(300, 11).Value2 = 12000
$ws.Cells.Item(300, 12).Value2 = 13000
$ws.Cells.Item(300, 13).Value2 = 12500
$ws.Cells.Item(300, 16).Value2 = 125
